$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update sub-table of section "1 Planejamento" (rows 4-8) ---
# The task "1.2 Criação do Documento de visão do Projeto" (row 5) is removed;
# subsequent tasks (1.3, 1.4, 1.5) shift up one position (to 1.2, 1.3, 1.4)
# and the last row becomes empty (but keeps its border/format).

# 1.1 task description is renamed
$ws.Range("C4").Value2 = "Escopo do Produto e Projeto"

# Capture the values that need to move up before overwriting them
$d6 = $ws.Range("D6").Value2
$d7 = $ws.Range("D7").Value2
$d8 = $ws.Range("D8").Value2

# Row 5 (was 1.2 "Criação do Documento de visão do Projeto") becomes 1.2 "Identificação de Requisitos de Negócios e Técnicos"
$ws.Range("C5").Value2 = "Identificação de Requisitos de Negócios e Técnicos"
$ws.Range("D5").Value2 = $d6

# Row 6 (was 1.3 "Identificação de Requisitos de Negócios e Técnicos") becomes 1.3 "Desenvolvimento do Plano de Projeto"
$ws.Range("C6").Value2 = "Desenvolvimento do Plano de Projeto"
$ws.Range("D6").Value2 = $d7

# Row 7 (was 1.4 "Desenvolvimento do Plano de Projeto") becomes 1.4 "Identificação dos Recursos necessários"
$ws.Range("C7").Value2 = "Identificação dos Recursos necessários"
$ws.Range("D7").Value2 = $d8
# Row 7 takes on the "last row of table" formatting (thicker bottom border), copied from row 8
$ws.Range("B8:D8").Copy() | Out-Null
$ws.Range("B7:D7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# Row 8 (was 1.5 "Identificação dos Recursos necessários") is now empty
$ws.Range("B8").Value2 = $null
$ws.Range("C8").Value2 = $null
$ws.Range("D8").Value2 = $null

# --- Update sheet view (scroll position / active cell selection) ---
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F8").Select() | Out-Null

$wb.Save()
